# Applies the cryptos.xlsx price/volume refresh described in the commit message:
# 'Updated cryptos list on Fri Jun 23 18:40:19 UTC 2023 with GitHub Actions'
#
# Column D/E cells are stored as text in the workbook (t="inlineStr"), and some
# of the new values (e.g. "1.003") look numeric, so a leading apostrophe is
# concatenated onto the value to force Excel to keep storing them as text,
# exactly like the original data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '31.003.87'
$ws.Range('E2').Value = "'" + '  +3.18%  '
$ws.Range('D3').Value = "'" + '1.909.49'
$ws.Range('E3').Value = "'" + '  +1.31%  '
$ws.Range('D4').Value = "'" + '1.003'
$ws.Range('E4').Value = "'" + '  +0.47%  '
$ws.Range('D5').Value = "'" + '245.64'
$ws.Range('E5').Value = "'" + '  +0.49%  '
$ws.Range('D7').Value = "'" + '0.4981'
$ws.Range('E7').Value = "'" + '  +0.35%  '
$ws.Range('D8').Value = "'" + '0.2992'
$ws.Range('E8').Value = "'" + '  +2.31%  '
$ws.Range('D9').Value = "'" + '0.06874'
$ws.Range('E9').Value = "'" + '  +3.69%  '
$ws.Range('D10').Value = "'" + '1.913.66'
$ws.Range('E10').Value = "'" + '  +1.58%  '
$ws.Range('D11').Value = "'" + '16.97'
$ws.Range('E11').Value = "'" + '  -0.25%  '
$ws.Range('D12').Value = "'" + '0.07305'
$ws.Range('E12').Value = "'" + '  +1.41%  '
$ws.Range('D13').Value = "'" + '91.02'
$ws.Range('E13').Value = "'" + '  +5.78%  '
$ws.Range('D14').Value = "'" + '5.087'
$ws.Range('E14').Value = "'" + '  +4.69%  '
$ws.Range('D15').Value = "'" + '0.6797'
$ws.Range('E15').Value = "'" + '  +1.69%  '
$ws.Range('D16').Value = "'" + '30.997.39'
$ws.Range('E16').Value = "'" + '  +3.24%  '
$ws.Range('D17').Value = "'" + '0.000008039'
$ws.Range('E17').Value = "'" + '  +2.49%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').Value = "'" + '1.002'
$ws.Range('E18').Value = "'" + '  +0.35%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').Value = "'" + '13.28'
$ws.Range('E19').Value = "'" + '  +3.32%  '
$ws.Range('D20').Value = "'" + '2.163.95'
$ws.Range('E20').Value = "'" + '  +2.00%  '
$ws.Range('D21').Value = "'" + '0.9996'
$ws.Range('E21').Value = "'" + '  +0.23%  '
$ws.Range('D22').Value = "'" + '4.867'
$ws.Range('E22').Value = "'" + '  +1.89%  '
$ws.Range('D23').Value = "'" + '183.31'
$ws.Range('E23').Value = "'" + '  +34.38%  '
$ws.Range('D24').Value = "'" + '6.080'
$ws.Range('E24').Value = "'" + '  +8.27%  '
$ws.Range('D25').Value = "'" + '9.360'
$ws.Range('E25').Value = "'" + '  +1.82%  '
$ws.Range('D26').Value = "'" + '153.43'
$ws.Range('E26').Value = "'" + '  +2.42%  '
$ws.Range('D27').Value = "'" + '18.75'
$ws.Range('E27').Value = "'" + '  +11.69%  '
$ws.Range('D28').Value = "'" + '1.942'
$ws.Range('E28').Value = "'" + '  +1.60%  '
$ws.Range('D29').Value = "'" + '1.404'
$ws.Range('E29').Value = "'" + '  +1.66%  '
$ws.Range('D30').Value = "'" + '4.342'
$ws.Range('E30').Value = "'" + '  +3.60%  '
$ws.Range('D31').Value = "'" + '0.08963'
$ws.Range('E31').Value = "'" + '  +3.27%  '
$ws.Range('D32').Value = "'" + '4.049'
$ws.Range('E32').Value = "'" + '  +2.09%  '
$ws.Range('D33').Value = "'" + '0.05245'
$ws.Range('E33').Value = "'" + '  +4.97%  '
$ws.Range('D34').Value = "'" + '0.7479'
$ws.Range('E34').Value = "'" + '  +6.19%  '
$ws.Range('D35').Value = "'" + '1.144'
$ws.Range('E35').Value = "'" + '  +3.17%  '
$ws.Range('D36').Value = "'" + '2.673'
$ws.Range('E36').Value = "'" + '  +0.78%  '
$ws.Range('D37').Value = "'" + '0.01935'
$ws.Range('E37').Value = "'" + '  +17.54%  '
$ws.Range('D38').Value = "'" + '2.732'
$ws.Range('E38').Value = "'" + '  +1.33%  '
$ws.Range('D39').Value = "'" + '2.181'
$ws.Range('E39').Value = "'" + '  -1.19%  '
$ws.Range('D40').Value = "'" + '0.9366'
$ws.Range('E40').Value = "'" + '  +0.13%  '
$ws.Range('D41').Value = "'" + '0.4371'
$ws.Range('E41').Value = "'" + '  +4.24%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = "'" + '106.02'
$ws.Range('E42').Value = "'" + '  +4.22%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = "'" + '5.875'
$ws.Range('E43').Value = "'" + '  -1.61%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = "'" + '1.002'
$ws.Range('E44').Value = "'" + '  +0.22%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = "'" + '7.792'
$ws.Range('E45').Value = "'" + '  +2.70%  '
$ws.Range('D46').Value = "'" + '0.1343'
$ws.Range('E46').Value = "'" + '  +6.19%  '
$ws.Range('D47').Value = "'" + '0.05860'
$ws.Range('E47').Value = "'" + '  +2.33%  '
$ws.Range('D48').Value = "'" + '8.616'
$ws.Range('E48').Value = "'" + '  +4.78%  '
$ws.Range('D49').Value = "'" + '0.3889'
$ws.Range('E49').Value = "'" + '  +4.84%  '
$ws.Range('D50').Value = "'" + '33.27'
$ws.Range('E50').Value = "'" + '  +2.46%  '
$ws.Range('D51').Value = "'" + '1.391'
$ws.Range('E51').Value = "'" + '  +4.07%  '
